$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "694÷9=" "439÷7="
Replace-Text "810÷5=" "750÷4="
Replace-Text "578÷9=" "882÷3="
Replace-Text "217÷4=" "602÷8="
Replace-Text "469÷5=" "139÷4="
Replace-Text "541÷3=" "810÷5="
Replace-Text "126÷8=" "369÷7="
Replace-Text "139÷3=" "252÷6="
Replace-Text "506÷4=" "222÷3="
Replace-Text "219÷2=" "617÷7="
Replace-Text "772÷4=" "838÷8="
Replace-Text "121÷3=" "278÷5="
Replace-Text "210÷2=" "581÷6="
Replace-Text "590÷7=" "890÷2="
Replace-Text "555÷6=" "959÷6="
Replace-Text "202÷4=" "459÷3="
Replace-Text "626÷2=" "716÷8="
Replace-Text "144÷7=" "354÷8="
Replace-Text "324÷4=" "831÷4="
Replace-Text "262÷6=" "545÷9="
Replace-Text "696÷4=" "278÷3="
Replace-Text "962÷7=" "495÷6="
Replace-Text "688÷4=" "363÷5="
Replace-Text "844÷5=" "535÷7="
Replace-Text "113÷7=" "600÷4="
